$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")

# --- Add new row 86: UK survey round 73, panel F, wave 27 ---
$ws.Range("A86").Value2 = 3
$ws.Range("B86").Value2 = 0
$ws.Range("C86").Value2 = "uk"
$ws.Range("D86").Value2 = 73
$ws.Range("E86").Value2 = "F"
$ws.Range("F86").Value2 = 27
$ws.Range("G86").Value2 = 44427
$ws.Range("G85").Copy()
$ws.Range("G86").PasteSpecial(-4122)
$ws.Range("H86").Value2 = "21-037554_PFW27_Final_ICUO"
$ws.Range("I86").Formula = '=C86&"_"&"sr"&TEXT(D86,"00")&"_"&YEAR(G86)&TEXT(G86,"MM")&TEXT(G86,"DD")&"_p"&E86&"_wv"&TEXT(F86,"00")&""'
$ws.Range("J86").Value2 = 1

# --- Add new row 87: UK survey round 74, panel E, wave 28 ---
$ws.Range("A87").Value2 = 3
$ws.Range("B87").Value2 = 0
$ws.Range("C87").Value2 = "uk"
$ws.Range("D87").Value2 = 74
$ws.Range("E87").Value2 = "E"
$ws.Range("F87").Value2 = 28
$ws.Range("G87").Value2 = 44434
$ws.Range("G85").Copy()
$ws.Range("G87").PasteSpecial(-4122)
$ws.Range("H87").Value2 = "21-037558_PEW28_Final_ICUO"
$ws.Range("I87").Formula = '=C87&"_"&"sr"&TEXT(D87,"00")&"_"&YEAR(G87)&TEXT(G87,"MM")&TEXT(G87,"DD")&"_p"&E87&"_wv"&TEXT(F87,"00")&""'

$excel.CutCopyMode = 0

# --- Update the view so the newly added rows are visible/selected ---
$ws.Range("J89").Select()
